$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: was "level 4" -> becomes "strength 1" (values shifted/reduced) ---
$ws.Range("A5:B5").Clear()
$ws.Range("C5").Value = "strength 1"
$ws.Range("D5").Value = 19
$ws.Range("G5").Value = 11
$ws.Range("J5").Value = 16

# --- Row 6: was "strength 1" -> becomes "strength 2" (values shifted/reduced) ---
$ws.Range("C6").Value = "strength 2"
$ws.Range("D6").Value = 20
$ws.Range("G6").Value = 12
$ws.Range("J6").Value = 17
$ws.Range("P6").Value = 15

# --- Row 7: was "strength 2" -> becomes "abyss" (greatsword/fury/blade_of_dread columns removed) ---
$ws.Range("C7").Value = "abyss"
$ws.Range("D7:L7").Clear()
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = 3
$ws.Range("P7").Value = 12
$ws.Range("Q7").Value = 4

# --- Row 8: was "abyss 1" -> becomes "darkness" (tenebris_touch columns removed) ---
$ws.Range("C8").Value = "darkness"
$ws.Range("M8:O8").Clear()

# --- Rows 9-11 ("abyss 2", "darkness 1", "darkness 2") are removed entirely ---
$ws.Rows("9:11").Delete()

# --- Restore a plain selection state ---
$ws.Range("Q9").Select() | Out-Null
